# "adding averages and more checks"
#
# - Training Dashboard: LAST UPDATE (col I) rolls forward 8 days,
#   08-Sep-2025 -> 16-Sep-2025, for every data row (3-19); PERIOD TO
#   EXPIRE (col H) drops by 8 for each of those rows to match.
# - Exam Dashboard: E3:E5 comments reworded from "OK" to "date is valid",
#   and column E is widened to fit the longer text.
# - Header styling: title + table-header rows switch to bold white text
#   (drops the old 14pt title size in favour of bold+white everywhere).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Training Dashboard
$ws2 = $wb.Worksheets.Item(2)   # Exam Dashboard

# ---------------------------------------------------------------------
# Training Dashboard: rows 3-19, columns H (PERIOD TO EXPIRE) and
# I (LAST UPDATE)
# ---------------------------------------------------------------------
$periodToExpire = @(590, 591, 594, 594, 590, 590, 591, 594, 710, 591, 706, -19618, 278, 314, 314, 308, 348)

# I already holds plain text dates ("08-Sep-2025"); keep it text so the
# new date isn't silently reinterpreted as a serial date number.
$ws1.Range("I3:I19").NumberFormat = "@"

for ($i = 0; $i -lt $periodToExpire.Length; $i++) {
    $row = $i + 3
    $ws1.Cells.Item($row, 8).Value = $periodToExpire[$i]
    $ws1.Cells.Item($row, 9).Value = "16-Sep-2025"
}

# ---------------------------------------------------------------------
# Exam Dashboard: reworded remarks + wider column E
# ---------------------------------------------------------------------
$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"
$ws2.Range("E5").Value = "date is valid"

$ws2.Range("E1").ColumnWidth = 14.17

# ---------------------------------------------------------------------
# Styles: title (A1) and header rows (row 2) become bold white text
# ---------------------------------------------------------------------
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
$ws1.Range("A2:K2").Font.Color = 16777215

$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

Write-Host "Applied: refreshed LAST UPDATE/PERIOD TO EXPIRE, reworded exam remarks, widened column E, restyled headers."
